$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI values following Dr Hou's advice (Ligand-expressing / Receptor-expressing
# cell counts changed from 1 to 3, with all derived statistics recomputed accordingly).

$data = @{
    2  = @{ E=3; G=141.3574043333333;  H=424.072213;         I=0.4954750229273862; J=0.4954750229273862; K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636; Q=917.4869708617733;  R=8257.382737755959; S=0.3974535504295279;  T=0.3974535504295278 }
    3  = @{ E=3; G=141.3574043333333;  H=424.072213;         I=0.4954750229273862; J=0.4954750229273862; K=3; M=1.600723;          N=4.802169;  O=0.1978333275383364; P=0.1978333275383364; Q=226.2740483366663;  R=2036.466435029997; S=0.09802147249785832; T=0.09802147249785831 }
    4  = @{ E=3; G=140.4344916666667;  H=421.303475;         I=0.4922401009448182; J=0.4922401009448182; K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636; Q=911.4967622066029;  R=8203.470859859424; S=0.3948586038270983;  T=0.3948586038270983 }
    5  = @{ E=3; G=140.4344916666667;  H=421.303475;         I=0.4922401009448182; J=0.4922401009448182; K=3; M=1.600723;          N=4.802169;  O=0.1978333275383364; P=0.1978333275383364; Q=224.7967208041417;  R=2023.170487237275; S=0.09738149711772;    T=0.09738149711771998 }
    6  = @{ E=3; G=0.248853;           H=0.7465590000000001; I=0.0008722602573388757; J=0.0008722602573388757; K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636; Q=1.615192258493;    R=14.536730326437;   S=0.0006996981081500805; T=0.0006996981081500804 }
    7  = @{ E=3; G=0.248853;           H=0.7465590000000001; I=0.0008722602573388757; J=0.0008722602573388757; K=3; M=1.600723;          N=4.802169;  O=0.1978333275383364; P=0.1978333275383364; Q=0.3983447207190001; R=3.585102486471;    S=0.0001725621491887954; T=0.0001725621491887954 }
    8  = @{ E=3; G=3.255981999999999;  H=9.767945999999998;  I=0.01141261587045664; J=0.01141261587045664; K=3; M=6.490547666666667; N=19.471643; O=0.8021666724616637; P=0.8021666724616636; Q=21.13310637280867;  R=190.197957355278;  S=0.009154820096887378; T=0.009154820096887378 }
    9  = @{ E=3; G=3.255981999999999;  H=9.767945999999998;  I=0.01141261587045664; J=0.01141261587045664; K=3; M=1.600723;          N=4.802169;  O=0.1978333275383364; P=0.1978333275383364; Q=5.211925274985999;  R=46.90732747487399; S=0.002257795773569265; T=0.002257795773569265 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}
